$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is being extended with a new column R holding the 2021 data
# point (mirroring the existing 2007-2020 columns D..Q). Copy the
# formatting of column Q (rows 4-34, the data block with its header,
# borders, number format, etc.) into the brand-new column R so every
# new cell inherits the same look as its neighbour in the same row.
$ws.Range("Q4:Q34").Copy()
$ws.Range("R4:R34").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Header year for the new column.
$ws.Cells.Item(4, 18).Value = 2021

# Data values for the new column, one per data row.
$ws.Cells.Item(5, 18).Value = 11.9
$ws.Cells.Item(6, 18).Value = 13.1
$ws.Cells.Item(7, 18).Value = 10.6
$ws.Cells.Item(8, 18).Value = 11
$ws.Cells.Item(9, 18).Value = 10
$ws.Cells.Item(10, 18).Value = 12
$ws.Cells.Item(11, 18).Value = 10.199999999999999
$ws.Cells.Item(12, 18).Value = 10.5
$ws.Cells.Item(13, 18).Value = 10
$ws.Cells.Item(14, 18).Value = 19.399999999999999
$ws.Cells.Item(15, 18).Value = 22.3
$ws.Cells.Item(16, 18).Value = 16.399999999999999
$ws.Cells.Item(17, 18).Value = 9.4
$ws.Cells.Item(18, 18).Value = 11.4
$ws.Cells.Item(19, 18).Value = 7.3
$ws.Cells.Item(20, 18).Value = 3.1
$ws.Cells.Item(21, 18).Value = 2.9
$ws.Cells.Item(22, 18).Value = 3.4
$ws.Cells.Item(23, 18).Value = 15
$ws.Cells.Item(24, 18).Value = 17.3
$ws.Cells.Item(25, 18).Value = 12.7
$ws.Cells.Item(26, 18).Value = 7.9
$ws.Cells.Item(27, 18).Value = 8.4
$ws.Cells.Item(28, 18).Value = 7.4
$ws.Cells.Item(29, 18).Value = 15.2
$ws.Cells.Item(30, 18).Value = 17.600000000000001
$ws.Cells.Item(31, 18).Value = 12.6
$ws.Cells.Item(32, 18).Value = 27.9
$ws.Cells.Item(33, 18).Value = 32.700000000000003
$ws.Cells.Item(34, 18).Value = 22.8

# Reflect the recorded cursor position after the edit (cell R3 selected).
$ws.Range("R3").Select()
